$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totals row (9): swap manual +-chains for AutoSum-style SUM() calls ---
$ws.Range("B9").Formula = "=SUM(B4:B8)"
$ws.Range("C9:D9").Formula = "=SUM(C4:C8)"

# --- Column E (row total): AutoSum filled down E5:E9 in one go, turning
#     it into a shared formula that also re-derives E9 ---
$ws.Range("E5:E9").Formula = "=SUM(B5:D5)"

# --- New rows: MIN / MAX / AVERAGE / COUNT summary block ---
# Each column formula is entered independently (not filled as one block),
# so every cell keeps its own standalone formula rather than a shared one.
$ws.Range("A11").Value = "MIN"
$ws.Range("B11").Formula = "=MIN(B4:B8)"
$ws.Range("C11").Formula = "=MIN(C4:C8)"
$ws.Range("D11").Formula = "=MIN(D4:D8)"
$ws.Range("E11").Formula = "=MIN(E4:E8)"

$ws.Range("A12").Value = "MAX"
$ws.Range("B12").Formula = "=MAX(B4:B8)"
$ws.Range("C12").Formula = "=MAX(C4:C8)"
$ws.Range("D12").Formula = "=MAX(D4:D8)"
$ws.Range("E12").Formula = "=MAX(E4:E8)"

$ws.Range("A13").Value = "AVERAGE"
$ws.Range("B13").Formula = "=AVERAGE(B4:B8)"
$ws.Range("C13").Formula = "=AVERAGE(C4:C8)"
$ws.Range("D13").Formula = "=AVERAGE(D4:D8)"
$ws.Range("E13").Formula = "=AVERAGE(E4:E8)"

$ws.Range("A14").Value = "COUNT"
$ws.Range("B14").Formula = "=COUNT(B4:B8)"
$ws.Range("C14").Formula = "=COUNT(C4:C8)"
$ws.Range("D14").Formula = "=COUNT(D4:D8)"
$ws.Range("E14").Formula = "=COUNT(E4:E8)"

# give F11 the percent style, matching the rest of column F
$ws.Range("F11").NumberFormat = $ws.Range("F9").NumberFormat

# update selection to reflect where the user left off
$ws.Range("B18").Select() | Out-Null
